$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 29 (RelAcctName): shorten Chinese name, move old text to remarks column
$ws.Range("C29").Value = "帳戶戶名"
$ws.Range("G29").Value = "第三人帳戶戶名"

# Row 30 (RelationId): shorten Chinese name, move old text to remarks column
$ws.Range("C30").Value = "身分證字號"
$ws.Range("G30").Value = "第三人身分證字號"

# Row 31 (RelAcctBirthday): shorten Chinese name, move old text to remarks column
$ws.Range("C31").Value = "出生日期"
$ws.Range("G31").Value = "第三人出生日期"

# Row 32 (RelAcctGender): shorten Chinese name, prepend old text to remarks column
$ws.Range("C32").Value = "性別"
$ws.Range("G32").Value = "第三人性別" + [char]10 + "CdCode.Sex"
$ws.Range("G32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 32.4

# Restore view state: scroll position and active selection
$ws.Range("C31").Select()
$excel.ActiveWindow.ScrollRow = 26
